$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.46910000000002
$ws.Range("C4").Value = -11.08089999999999
$ws.Range("E4").Value = 13.0566

$ws.Range("C5").Value = -14.63690000000001

$ws.Range("A7").Value = -21.44070000000001

$ws.Range("C8").Value = -11.7557

$ws.Range("E9").Value = 13.79770000000001

$ws.Range("A16").Value = -20.24589999999998
$ws.Range("C16").Value = -12.0436

$ws.Range("E18").Value = 13.0821
